$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: 5, 2, 5, chaussure
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = "chaussure"

# Add new row 8: 7, 15, 15, radis
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 15
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = "radis"

# Update selection to match the diff
$ws.Range("E16").Select()
